$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Cases by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value = 270
$ws.Range("B3").Value = 1311
$ws.Range("B4").Value = 3583
$ws.Range("B5").Value = 15423
$ws.Range("B6").Value = 16974
$ws.Range("B7").Value = 14882
$ws.Range("B8").Value = 12496
$ws.Range("B9").Value = 4496
$ws.Range("B10").Value = 3032
$ws.Range("B11").Value = 1817
$ws.Range("B12").Value = 1186
$ws.Range("B13").Value = 1867
[void]$ws.Range("B21").Select()

# ---------------------------------------------------------------------------
# Sheet: Cases by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 26035
$ws.Range("B3").Value = 50417
$ws.Range("B4").Value = 898
[void]$ws.Range("B16").Select()

# ---------------------------------------------------------------------------
# Sheet: Cases by RaceEthnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 950
$ws.Range("B3").Value = 12893
$ws.Range("B4").Value = 28009
$ws.Range("B5").Value = 491
$ws.Range("B6").Value = 26507
$ws.Range("B7").Value = 8500
[void]$ws.Range("C20").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B3").Value = 12
$ws.Range("B4").Value = 33
$ws.Range("B5").Value = 233
$ws.Range("B6").Value = 785
$ws.Range("B7").Value = 2305
$ws.Range("B8").Value = 5321
$ws.Range("B9").Value = 4488
$ws.Range("B10").Value = 5818
$ws.Range("B11").Value = 6455
$ws.Range("B12").Value = 6414
$ws.Range("B13").Value = 16239
[void]$ws.Range("B22").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 20162
$ws.Range("B3").Value = 27946
[void]$ws.Range("C16").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1006
$ws.Range("B3").Value = 4796
$ws.Range("B4").Value = 22329
$ws.Range("B5").Value = 262
$ws.Range("B6").Value = 19691
$ws.Activate()
[void]$ws.Range("D22").Select()
